$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date for all data rows
$updatedRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,111,112,113,114,115,116,117,118,119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,134,135,136,137,138,139,140,141)
foreach ($r in $updatedRows) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# Add display-text argument to HYPERLINK formulas in columns S,T,V,W,X,Y for rows 2-21
# Row 2: A 33217-2020
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 33217-2020.xlsx", "A 33217-2020")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 33217-2020.png", "A 33217-2020")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 33217-2020.docx", "A 33217-2020")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 33217-2020.docx", "A 33217-2020")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 33217-2020.docx", "A 33217-2020")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 33217-2020.docx", "A 33217-2020")'

# Row 3: A 12071-2023
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 12071-2023.xlsx", "A 12071-2023")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 12071-2023.png", "A 12071-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 12071-2023.docx", "A 12071-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 12071-2023.docx", "A 12071-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 12071-2023.docx", "A 12071-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 12071-2023.docx", "A 12071-2023")'

# Row 4: A 26954-2022
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 26954-2022.xlsx", "A 26954-2022")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 26954-2022.png", "A 26954-2022")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 26954-2022.docx", "A 26954-2022")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 26954-2022.docx", "A 26954-2022")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 26954-2022.docx", "A 26954-2022")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 26954-2022.docx", "A 26954-2022")'

# Row 5: A 59192-2021
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 59192-2021.xlsx", "A 59192-2021")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 59192-2021.png", "A 59192-2021")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 59192-2021.docx", "A 59192-2021")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 59192-2021.docx", "A 59192-2021")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 59192-2021.docx", "A 59192-2021")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 59192-2021.docx", "A 59192-2021")'

# Row 6: A 62068-2019
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 62068-2019.xlsx", "A 62068-2019")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 62068-2019.png", "A 62068-2019")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 62068-2019.docx", "A 62068-2019")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 62068-2019.docx", "A 62068-2019")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 62068-2019.docx", "A 62068-2019")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 62068-2019.docx", "A 62068-2019")'

# Row 7: A 60734-2020
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 60734-2020.xlsx", "A 60734-2020")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 60734-2020.png", "A 60734-2020")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 60734-2020.docx", "A 60734-2020")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 60734-2020.docx", "A 60734-2020")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 60734-2020.docx", "A 60734-2020")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 60734-2020.docx", "A 60734-2020")'

# Row 8: A 59199-2021
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 59199-2021.xlsx", "A 59199-2021")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 59199-2021.png", "A 59199-2021")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 59199-2021.docx", "A 59199-2021")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 59199-2021.docx", "A 59199-2021")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 59199-2021.docx", "A 59199-2021")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 59199-2021.docx", "A 59199-2021")'

# Row 9: A 25229-2023
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 25229-2023.xlsx", "A 25229-2023")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 25229-2023.png", "A 25229-2023")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 25229-2023.docx", "A 25229-2023")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 25229-2023.docx", "A 25229-2023")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 25229-2023.docx", "A 25229-2023")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 25229-2023.docx", "A 25229-2023")'

# Row 10: A 47035-2022
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 47035-2022.xlsx", "A 47035-2022")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 47035-2022.png", "A 47035-2022")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 47035-2022.docx", "A 47035-2022")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 47035-2022.docx", "A 47035-2022")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 47035-2022.docx", "A 47035-2022")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 47035-2022.docx", "A 47035-2022")'

# Row 11: A 61660-2018
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 61660-2018.xlsx", "A 61660-2018")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 61660-2018.png", "A 61660-2018")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 61660-2018.docx", "A 61660-2018")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 61660-2018.docx", "A 61660-2018")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 61660-2018.docx", "A 61660-2018")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 61660-2018.docx", "A 61660-2018")'

# Row 12: A 41183-2022
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 41183-2022.xlsx", "A 41183-2022")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 41183-2022.png", "A 41183-2022")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 41183-2022.docx", "A 41183-2022")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 41183-2022.docx", "A 41183-2022")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 41183-2022.docx", "A 41183-2022")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 41183-2022.docx", "A 41183-2022")'

# Row 13: A 64166-2021
$ws.Range("S13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 64166-2021.xlsx", "A 64166-2021")'
$ws.Range("T13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 64166-2021.png", "A 64166-2021")'
$ws.Range("V13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 64166-2021.docx", "A 64166-2021")'
$ws.Range("W13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 64166-2021.docx", "A 64166-2021")'
$ws.Range("X13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 64166-2021.docx", "A 64166-2021")'
$ws.Range("Y13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 64166-2021.docx", "A 64166-2021")'

# Row 14: A 15171-2022
$ws.Range("S14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 15171-2022.xlsx", "A 15171-2022")'
$ws.Range("T14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 15171-2022.png", "A 15171-2022")'
$ws.Range("V14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 15171-2022.docx", "A 15171-2022")'
$ws.Range("W14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 15171-2022.docx", "A 15171-2022")'
$ws.Range("X14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 15171-2022.docx", "A 15171-2022")'
$ws.Range("Y14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 15171-2022.docx", "A 15171-2022")'

# Row 15: A 35512-2018
$ws.Range("S15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 35512-2018.xlsx", "A 35512-2018")'
$ws.Range("T15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 35512-2018.png", "A 35512-2018")'
$ws.Range("V15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 35512-2018.docx", "A 35512-2018")'
$ws.Range("W15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 35512-2018.docx", "A 35512-2018")'
$ws.Range("X15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 35512-2018.docx", "A 35512-2018")'
$ws.Range("Y15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 35512-2018.docx", "A 35512-2018")'

# Row 16: A 45095-2019
$ws.Range("S16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 45095-2019.xlsx", "A 45095-2019")'
$ws.Range("T16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 45095-2019.png", "A 45095-2019")'
$ws.Range("V16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 45095-2019.docx", "A 45095-2019")'
$ws.Range("W16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 45095-2019.docx", "A 45095-2019")'
$ws.Range("X16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 45095-2019.docx", "A 45095-2019")'
$ws.Range("Y16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 45095-2019.docx", "A 45095-2019")'

# Row 17: A 14684-2021
$ws.Range("S17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 14684-2021.xlsx", "A 14684-2021")'
$ws.Range("T17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 14684-2021.png", "A 14684-2021")'
$ws.Range("V17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 14684-2021.docx", "A 14684-2021")'
$ws.Range("W17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 14684-2021.docx", "A 14684-2021")'
$ws.Range("X17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 14684-2021.docx", "A 14684-2021")'
$ws.Range("Y17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 14684-2021.docx", "A 14684-2021")'

# Row 18: A 17549-2023
$ws.Range("S18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 17549-2023.xlsx", "A 17549-2023")'
$ws.Range("T18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 17549-2023.png", "A 17549-2023")'
$ws.Range("V18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 17549-2023.docx", "A 17549-2023")'
$ws.Range("W18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 17549-2023.docx", "A 17549-2023")'
$ws.Range("X18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 17549-2023.docx", "A 17549-2023")'
$ws.Range("Y18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 17549-2023.docx", "A 17549-2023")'

# Row 19: A 60735-2018
$ws.Range("S19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 60735-2018.xlsx", "A 60735-2018")'
$ws.Range("T19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 60735-2018.png", "A 60735-2018")'
$ws.Range("V19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 60735-2018.docx", "A 60735-2018")'
$ws.Range("W19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 60735-2018.docx", "A 60735-2018")'
$ws.Range("X19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 60735-2018.docx", "A 60735-2018")'
$ws.Range("Y19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 60735-2018.docx", "A 60735-2018")'

# Row 20: A 45324-2022
$ws.Range("S20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 45324-2022.xlsx", "A 45324-2022")'
$ws.Range("T20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 45324-2022.png", "A 45324-2022")'
$ws.Range("V20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 45324-2022.docx", "A 45324-2022")'
$ws.Range("W20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 45324-2022.docx", "A 45324-2022")'
$ws.Range("X20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 45324-2022.docx", "A 45324-2022")'
$ws.Range("Y20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 45324-2022.docx", "A 45324-2022")'

# Row 21: A 25279-2023
$ws.Range("S21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/artfynd/A 25279-2023.xlsx", "A 25279-2023")'
$ws.Range("T21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/kartor/A 25279-2023.png", "A 25279-2023")'
$ws.Range("V21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomål/A 25279-2023.docx", "A 25279-2023")'
$ws.Range("W21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/klagomålsmail/A 25279-2023.docx", "A 25279-2023")'
$ws.Range("X21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsyn/A 25279-2023.docx", "A 25279-2023")'
$ws.Range("Y21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ARJEPLOG/tillsynsmail/A 25279-2023.docx", "A 25279-2023")'

